$d = $word.ActiveDocument
$full = $d.Content.Text

# --- Change 1: "For each" highlight yellow -> red (first occurrence) ---
$idx = $full.IndexOf("For each")
$r = $d.Range($idx, $idx + "For each".Length)
$r.Font.HighlightColorIndex = 6  # wdRed

# --- Change 2: split " without the system. I" + "f the payment is invalid"
#     into " without the system. " + "If the payment is invalid" (magenta) ---
$anchor = " without the system. If the payment is invalid"
$idx = $full.IndexOf($anchor)
$start = $idx + " without the system. ".Length
$end = $start + "If the payment is invalid".Length
$r = $d.Range($start, $end)
$r.Font.HighlightColorIndex = 5  # wdPink (magenta)

# --- Change 3: split " clicks the Cancel Button to cancel the rental form"
#     into " clicks the Cancel Button to " + "cancel the rental form" (cyan) ---
$anchor = " clicks the Cancel Button to cancel the rental form"
$idx = $full.IndexOf($anchor)
$start = $idx + " clicks the Cancel Button to ".Length
$end = $start + "cancel the rental form".Length
$r = $d.Range($start, $end)
$r.Font.HighlightColorIndex = 3  # wdTurquoise (cyan)

# --- Change 4: highlight cyan across "the COMPLETE button" + " " + "to save"
#     (splits ", the clerk clicks the COMPLETE button" and "to save it in the system. T") ---
$anchor = ", the clerk clicks the COMPLETE button"
$idx = $full.IndexOf($anchor)
$start = $idx + ", the clerk clicks ".Length

$anchor2 = "to save it in the system"
$idx2 = $full.IndexOf($anchor2)
$end = $idx2 + "to save".Length

$r = $d.Range($start, $end)
$r.Font.HighlightColorIndex = 3  # wdTurquoise (cyan)
